$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, pushing the existing rows 21-22 down to 22-23.
$ws.Rows.Item(21).Insert()

# Fill the newly inserted row 21 with this week's data point.
$ws.Cells.Item(21, 1).Value = 5
$ws.Cells.Item(21, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(21, 3).Value = "Maule"
$ws.Cells.Item(21, 4).Value = 44726
$ws.Cells.Item(21, 5).Value = 7
$ws.Cells.Item(21, 6).Value = 100112040
$ws.Cells.Item(21, 7).Value = "Cilantro"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 150
$ws.Cells.Item(21, 11).Value = 8000
$ws.Cells.Item(21, 12).Value = 8000
$ws.Cells.Item(21, 13).Value = 8000
$ws.Cells.Item(21, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(21, 15).Value = "Región del Maule"
$ws.Cells.Item(21, 16).Value = 222
$ws.Cells.Item(21, 17).Value = 36
$ws.Cells.Item(21, 18).Value = "Hortaliza"
